$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete now-unused rows 41-48 (data shrinks from 48 to 40 data rows)
$ws.Range("A41:C48").EntireRow.Delete() | Out-Null

# Update rows 2-40 with the new Entsoe IGCC netting flow data (Elnet added to portfolio)
$ws.Cells.Item(2,1).Value = 45859
$ws.Cells.Item(2,2).Value = 0.051
$ws.Cells.Item(2,3).Value = 23.877
$ws.Cells.Item(3,1).Value = 45859.01041666666
$ws.Cells.Item(3,2).Value = 0
$ws.Cells.Item(3,3).Value = 24.149
$ws.Cells.Item(4,1).Value = 45859.02083333334
$ws.Cells.Item(4,2).Value = 0
$ws.Cells.Item(4,3).Value = 34.899
$ws.Cells.Item(5,1).Value = 45859.03125
$ws.Cells.Item(5,2).Value = 0
$ws.Cells.Item(5,3).Value = 18.882
$ws.Cells.Item(6,1).Value = 45859.04166666666
$ws.Cells.Item(6,2).Value = 0
$ws.Cells.Item(6,3).Value = 26.851
$ws.Cells.Item(7,1).Value = 45859.05208333334
$ws.Cells.Item(7,2).Value = 0
$ws.Cells.Item(7,3).Value = 17.738
$ws.Cells.Item(8,1).Value = 45859.0625
$ws.Cells.Item(8,2).Value = 0
$ws.Cells.Item(8,3).Value = 14.032
$ws.Cells.Item(9,1).Value = 45859.07291666666
$ws.Cells.Item(9,2).Value = 0
$ws.Cells.Item(9,3).Value = 27
$ws.Cells.Item(10,1).Value = 45859.08333333334
$ws.Cells.Item(10,2).Value = 0
$ws.Cells.Item(10,3).Value = 41.786
$ws.Cells.Item(11,1).Value = 45859.09375
$ws.Cells.Item(11,2).Value = 0
$ws.Cells.Item(11,3).Value = 25.606
$ws.Cells.Item(12,1).Value = 45859.10416666666
$ws.Cells.Item(12,2).Value = 0
$ws.Cells.Item(12,3).Value = 8.887
$ws.Cells.Item(13,1).Value = 45859.11458333334
$ws.Cells.Item(13,2).Value = 0
$ws.Cells.Item(13,3).Value = 2.219
$ws.Cells.Item(14,1).Value = 45859.125
$ws.Cells.Item(14,2).Value = 0
$ws.Cells.Item(14,3).Value = 30.926
$ws.Cells.Item(15,1).Value = 45859.13541666666
$ws.Cells.Item(15,2).Value = 0
$ws.Cells.Item(15,3).Value = 32.078
$ws.Cells.Item(16,1).Value = 45859.14583333334
$ws.Cells.Item(16,2).Value = 0
$ws.Cells.Item(16,3).Value = 4.498
$ws.Cells.Item(17,1).Value = 45859.15625
$ws.Cells.Item(17,2).Value = 0
$ws.Cells.Item(17,3).Value = 0.223
$ws.Cells.Item(18,1).Value = 45859.16666666666
$ws.Cells.Item(18,2).Value = 0
$ws.Cells.Item(18,3).Value = 0.206
$ws.Cells.Item(19,1).Value = 45859.17708333334
$ws.Cells.Item(19,2).Value = 0
$ws.Cells.Item(19,3).Value = 0.92
$ws.Cells.Item(20,1).Value = 45859.1875
$ws.Cells.Item(20,2).Value = 0
$ws.Cells.Item(20,3).Value = 1.654
$ws.Cells.Item(21,1).Value = 45859.19791666666
$ws.Cells.Item(21,2).Value = 0
$ws.Cells.Item(21,3).Value = 6.549
$ws.Cells.Item(22,1).Value = 45859.20833333334
$ws.Cells.Item(22,2).Value = 0
$ws.Cells.Item(22,3).Value = 19.788
$ws.Cells.Item(23,1).Value = 45859.21875
$ws.Cells.Item(23,2).Value = 0
$ws.Cells.Item(23,3).Value = 21.908
$ws.Cells.Item(24,1).Value = 45859.22916666666
$ws.Cells.Item(24,2).Value = 0
$ws.Cells.Item(24,3).Value = 25.837
$ws.Cells.Item(25,1).Value = 45859.23958333334
$ws.Cells.Item(25,2).Value = 0.035
$ws.Cells.Item(25,3).Value = 10.489
$ws.Cells.Item(26,1).Value = 45859.25
$ws.Cells.Item(26,2).Value = 0
$ws.Cells.Item(26,3).Value = 30.69
$ws.Cells.Item(27,1).Value = 45859.26041666666
$ws.Cells.Item(27,2).Value = 0
$ws.Cells.Item(27,3).Value = 21.12
$ws.Cells.Item(28,1).Value = 45859.27083333334
$ws.Cells.Item(28,2).Value = 0
$ws.Cells.Item(28,3).Value = 19.688
$ws.Cells.Item(29,1).Value = 45859.28125
$ws.Cells.Item(29,2).Value = 0
$ws.Cells.Item(29,3).Value = 46.056
$ws.Cells.Item(30,1).Value = 45859.29166666666
$ws.Cells.Item(30,2).Value = 0
$ws.Cells.Item(30,3).Value = 40.267
$ws.Cells.Item(31,1).Value = 45859.30208333334
$ws.Cells.Item(31,2).Value = 0
$ws.Cells.Item(31,3).Value = 47.139
$ws.Cells.Item(32,1).Value = 45859.3125
$ws.Cells.Item(32,2).Value = 0
$ws.Cells.Item(32,3).Value = 28.497
$ws.Cells.Item(33,1).Value = 45859.32291666666
$ws.Cells.Item(33,2).Value = 0
$ws.Cells.Item(33,3).Value = 52.398
$ws.Cells.Item(34,1).Value = 45859.33333333334
$ws.Cells.Item(34,2).Value = 0
$ws.Cells.Item(34,3).Value = 37.104
$ws.Cells.Item(35,1).Value = 45859.34375
$ws.Cells.Item(35,2).Value = 0
$ws.Cells.Item(35,3).Value = 50.603
$ws.Cells.Item(36,1).Value = 45859.35416666666
$ws.Cells.Item(36,2).Value = 0
$ws.Cells.Item(36,3).Value = 10.39
$ws.Cells.Item(37,1).Value = 45859.36458333334
$ws.Cells.Item(37,2).Value = 0.278
$ws.Cells.Item(37,3).Value = 7.785
$ws.Cells.Item(38,1).Value = 45859.375
$ws.Cells.Item(38,2).Value = 0.986
$ws.Cells.Item(38,3).Value = 4.334
$ws.Cells.Item(39,1).Value = 45859.38541666666
$ws.Cells.Item(39,2).Value = 0.051
$ws.Cells.Item(39,3).Value = 7.415
$ws.Cells.Item(40,1).Value = 45859.39583333334
$ws.Cells.Item(40,2).Value = 0.241
$ws.Cells.Item(40,3).Value = 5.331
